$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 938
$ws.Range("B8").Value = "LC/GFG"
$ws.Range("C8").Value = "Count BST nodes that lie in a given range"

$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)

$ws.Range("C8").Select()
